$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume(1h) (E) columns to Text format while we
# write values, so numeric-looking strings (e.g. "1.000", "38.60") keep
# their exact textual representation instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '19.785.21'
$ws.Range('E2').Value = '  -8.67%  '
$ws.Range('D3').Value = '1.386.68'
$ws.Range('E3').Value = '  -9.65%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').Value = '268.04'
$ws.Range('E6').Value = '  -7.39%  '
$ws.Range('D7').Value = '0.3628'
$ws.Range('E7').Value = '  -7.82%  '
$ws.Range('D8').Value = '0.3038'
$ws.Range('E8').Value = '  -4.28%  '
$ws.Range('D9').Value = '39.24'
$ws.Range('E9').Value = '  -7.10%  '
$ws.Range('D10').Value = '0.9774'
$ws.Range('E10').Value = '  -7.40%  '
$ws.Range('D11').Value = '0.06421'
$ws.Range('E11').Value = '  -10.50%  '
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '5.295'
$ws.Range('E13').Value = '  -6.86%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').Value = '16.79'
$ws.Range('E14').Value = '  -9.55%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '6.065'
$ws.Range('E15').Value = '  -8.07%  '
$ws.Range('D16').Value = '1.388.01'
$ws.Range('E16').Value = '  -10.10%  '
$ws.Range('D17').Value = '0.000009939'
$ws.Range('E17').Value = '  -9.09%  '
$ws.Range('D18').Value = '0.05615'
$ws.Range('E18').Value = '  -15.00%  '
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = '69.80'
$ws.Range('E20').Value = '  -16.63%  '
$ws.Range('D21').Value = '5.505'
$ws.Range('E21').Value = '  -10.04%  '
$ws.Range('D22').Value = '14.42'
$ws.Range('E22').Value = '  -6.71%  '
$ws.Range('D23').Value = '10.55'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').Value = '2.238'
$ws.Range('E24').Value = '  -4.78%  '
$ws.Range('D25').Value = '19.784.79'
$ws.Range('E25').Value = '  -8.65%  '
$ws.Range('D26').Value = '2.176'
$ws.Range('E26').Value = '  -7.47%  '
$ws.Range('D27').Value = '136.69'
$ws.Range('E27').Value = '  -8.68%  '
$ws.Range('D28').Value = '16.53'
$ws.Range('E28').Value = '  -9.83%  '
$ws.Range('D29').Value = '1.541.72'
$ws.Range('E29').Value = '  -10.43%  '
$ws.Range('D30').Value = '107.61'
$ws.Range('E30').Value = '  -8.04%  '
$ws.Range('D31').Value = '3.825'
$ws.Range('E31').Value = '  -21.00%  '
$ws.Range('D32').Value = '5.217'
$ws.Range('E32').Value = '  -14.11%  '
$ws.Range('D33').Value = '0.7999'
$ws.Range('E33').Value = '  -14.91%  '
$ws.Range('D34').Value = '0.07575'
$ws.Range('D35').Value = '8.197'
$ws.Range('E35').Value = '  -3.81%  '
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').Value = '0.05625'
$ws.Range('E37').Value = '  -6.47%  '
$ws.Range('D38').Value = '4.680'
$ws.Range('E38').Value = '  -9.45%  '
$ws.Range('D39').Value = '0.02026'
$ws.Range('E39').Value = '  -9.01%  '
$ws.Range('D40').Value = '0.1869'
$ws.Range('E40').Value = '  -7.82%  '
$ws.Range('D41').Value = '10.00'
$ws.Range('E41').Value = '  -8.61%  '
$ws.Range('D42').Value = '1.295'
$ws.Range('E42').Value = '  -11.02%  '
$ws.Range('D43').Value = '1.048'
$ws.Range('E43').Value = '  -10.95%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.5185'
$ws.Range('E44').Value = '  -10.19%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '12.02'
$ws.Range('E45').Value = '  -7.58%  '
$ws.Range('D46').Value = '3.452'
$ws.Range('E46').Value = '  -7.04%  '
$ws.Range('D47').Value = '0.4980'
$ws.Range('E47').Value = '  -9.67%  '
$ws.Range('D48').Value = '109.07'
$ws.Range('E48').Value = '  -5.97%  '
$ws.Range('D49').Value = '1.725'
$ws.Range('E49').Value = '  -8.25%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '1.031'
$ws.Range('E50').Value = '  -11.57%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  +0.04%  '

# Restore the original (default) cell style now that the text values are
# stored, so formatting matches the source workbook.
$ws.Range("D2:E51").Style = "Normal"

